$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should match the formatting
# of the existing header cells (e.g. H1): bold font, border, centered.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for the new columns I (I0) and J (IF)
$data = @{
    2  = @(9, 9)
    3  = @(8, 8)
    4  = @(7, 7)
    5  = @(6, 8)
    6  = @(3, 5)
    7  = @(9, 9)
    8  = @(5, 5)
    9  = @(4, 5)
    10 = @(6, 7)
    11 = @(7, 8)
    12 = @(8, 8)
    13 = @(5, 5)
    14 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
